# Update Name of Algo
# Apply updated KNN imputation results to specific cells on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Value = -12.584
$ws.Range("A9").Value = -20.775
$ws.Range("C11").Value = -12.934
$ws.Range("A18").Value = -21.81
$ws.Range("A20").Value = -21.757
$ws.Range("D21").Value = -7.675999999999999
